# Auto-generated Excel COM-interop edit script
# Applies numeric updates to the Leve profit-calculation sheets
# (ALC, ARM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2164.16
$ws.Range("I40").Value = 1833.3334
$ws.Range("K40").Value = 1833.3334
$ws.Range("M40").Value = -1658.3334

$ws.Range("H62").Value = 110142.73
$ws.Range("I62").Value = 121056.4
$ws.Range("J62").Value = 1006
$ws.Range("K62").Value = 121056.4
$ws.Range("L62").Value = 1006
$ws.Range("M62").Value = -120432.4
$ws.Range("N62").Value = -2254

$ws.Range("H65").Value = 110142.73
$ws.Range("I65").Value = 121056.4
$ws.Range("J65").Value = 1006
$ws.Range("K65").Value = 605282
$ws.Range("L65").Value = 5030
$ws.Range("M65").Value = -602162
$ws.Range("N65").Value = -11270

$ws.Range("H106").Value = 2382.6
$ws.Range("I106").Value = 2249.1538
$ws.Range("K106").Value = 2249.1538
$ws.Range("M106").Value = -1618.1538

$ws.Range("H132").Value = 1885.2903
$ws.Range("I132").Value = 1928.8966
$ws.Range("J132").Value = 1253
$ws.Range("K132").Value = 5786.6898
$ws.Range("L132").Value = 3759
$ws.Range("M132").Value = -3256.6898
$ws.Range("N132").Value = -8819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -332

$ws.Range("H28").Value = 6232.6665
$ws.Range("I28").Value = 6232.6665
$ws.Range("K28").Value = 6232.6665
$ws.Range("M28").Value = -6040.6665

$ws.Range("H61").Value = 1874.75
$ws.Range("I61").Value = 750
$ws.Range("J61").Value = 2999.5
$ws.Range("K61").Value = 750
$ws.Range("L61").Value = 2999.5
$ws.Range("M61").Value = -538
$ws.Range("N61").Value = -3423.5

$ws.Range("H74").Value = 890.1818
$ws.Range("I74").Value = 931.0769
$ws.Range("J74").Value = 831.1111
$ws.Range("K74").Value = 931.0769
$ws.Range("L74").Value = 831.1111
$ws.Range("M74").Value = -57.07690000000002
$ws.Range("N74").Value = -2579.1111

$ws.Range("H77").Value = 890.1818
$ws.Range("I77").Value = 931.0769
$ws.Range("J77").Value = 831.1111
$ws.Range("K77").Value = 4655.3845
$ws.Range("L77").Value = 4155.555499999999
$ws.Range("M77").Value = -287.3845000000001
$ws.Range("N77").Value = -12891.5555

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

$ws.Range("H99").Value = 6232.6665
$ws.Range("I99").Value = 6232.6665
$ws.Range("K99").Value = 6232.6665
$ws.Range("M99").Value = -3237.6665

$ws.Range("H122").Value = 1798.7333
$ws.Range("I122").Value = 1613.9231
$ws.Range("K122").Value = 4841.7693
$ws.Range("M122").Value = -2391.7693

$ws.Range("H132").Value = 1373.5
$ws.Range("I132").Value = 883.3
$ws.Range("K132").Value = 2649.9
$ws.Range("M132").Value = -119.8999999999996

$ws.Range("H136").Value = 1874.75
$ws.Range("I136").Value = 750
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 2250
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = 300
$ws.Range("N136").Value = -14098.5

$ws.Range("H138").Value = 40000
$ws.Range("J138").Value = 40000
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2289.0908
$ws.Range("I105").Value = 2166.6667
$ws.Range("J105").Value = 2840
$ws.Range("K105").Value = 2166.6667
$ws.Range("L105").Value = 2840
$ws.Range("M105").Value = -419.6667000000002
$ws.Range("N105").Value = -6334

$ws.Range("H122").Value = 715178.4399999999
$ws.Range("I122").Value = 770115.25
$ws.Range("K122").Value = 2310345.75
$ws.Range("M122").Value = -2307895.75

$ws.Range("H132").Value = 3356.6667
$ws.Range("I132").Value = 2528
$ws.Range("K132").Value = 7584
$ws.Range("M132").Value = -5054

$ws.Range("H134").Value = 3050.4
$ws.Range("I134").Value = 2584
$ws.Range("J134").Value = 3750
$ws.Range("K134").Value = 7752
$ws.Range("L134").Value = 11250
$ws.Range("M134").Value = -5217
$ws.Range("N134").Value = -16320

$ws.Range("H140").Value = 53075
$ws.Range("J140").Value = 53075
$ws.Range("L140").Value = 53075
$ws.Range("N140").Value = -63435

$ws.Range("H141").Value = 45703.715
$ws.Range("J141").Value = 45703.715
$ws.Range("L141").Value = 45703.715
$ws.Range("N141").Value = -56063.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 913.37933
$ws.Range("I68").Value = 624.3333
$ws.Range("J68").Value = 1065.5088
$ws.Range("K68").Value = 1872.9999
$ws.Range("L68").Value = 3196.5264
$ws.Range("M68").Value = -1061.9999
$ws.Range("N68").Value = -4818.526400000001

$ws.Range("H71").Value = 913.37933
$ws.Range("I71").Value = 624.3333
$ws.Range("J71").Value = 1065.5088
$ws.Range("K71").Value = 5618.9997
$ws.Range("L71").Value = 9589.5792
$ws.Range("M71").Value = -1562.9997
$ws.Range("N71").Value = -17701.5792

$ws.Range("H86").Value = 376
$ws.Range("I86").Value = 376
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1128
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 58
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 376
$ws.Range("I89").Value = 376
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 3384
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 2544
$ws.Range("N89").ClearContents()

$ws.Range("H104").Value = 6000
$ws.Range("J104").Value = 6000
$ws.Range("L104").Value = 18000
$ws.Range("N104").Value = -23242

$ws.Range("H131").Value = 18074.672
$ws.Range("J131").Value = 1649.1666
$ws.Range("L131").Value = 4947.4998
$ws.Range("N131").Value = -15027.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4237.9062
$ws.Range("I70").Value = 4067.8
$ws.Range("J70").Value = 4845.4287
$ws.Range("K70").Value = 4067.8
$ws.Range("L70").Value = 4845.4287
$ws.Range("M70").Value = -3797.8
$ws.Range("N70").Value = -5385.4287

$ws.Range("H73").Value = 4237.9062
$ws.Range("I73").Value = 4067.8
$ws.Range("J73").Value = 4845.4287
$ws.Range("K73").Value = 4067.8
$ws.Range("L73").Value = 4845.4287
$ws.Range("M73").Value = -3131.8
$ws.Range("N73").Value = -6717.4287

$ws.Range("H132").Value = 1931.3478
$ws.Range("I132").Value = 1570.4054
$ws.Range("J132").Value = 3415.2222
$ws.Range("K132").Value = 4711.216200000001
$ws.Range("L132").Value = 10245.6666
$ws.Range("M132").Value = -2181.216200000001
$ws.Range("N132").Value = -15305.6666

$ws.Range("H138").Value = 19781.8
$ws.Range("J138").Value = 19781.8
$ws.Range("L138").Value = 19781.8
$ws.Range("N138").Value = -30061.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2476
$ws.Range("I93").Value = 2300
$ws.Range("K93").Value = 2300
$ws.Range("M93").Value = -1052

$ws.Range("H132").Value = 2978931
$ws.Range("I132").Value = 4904440.5
$ws.Range("J132").Value = 3143.4546
$ws.Range("K132").Value = 14713321.5
$ws.Range("L132").Value = 9430.363799999999
$ws.Range("M132").Value = -14710791.5
$ws.Range("N132").Value = -14490.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2222.6667
$ws.Range("I132").Value = 884.8333
$ws.Range("J132").Value = 3560.5
$ws.Range("K132").Value = 2654.4999
$ws.Range("L132").Value = 10681.5
$ws.Range("M132").Value = -124.4998999999998
$ws.Range("N132").Value = -15741.5
